$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999706558784,
    0.99753737692907074,
    0.98798131585940441,
    0.99098513655980836,
    0.98084297977282331,
    0.95687860249896428,
    0.95555549951327579,
    0.95342515465177535,
    0.95418489815809526,
    0.95590457001402385,
    0.95705445594230609,
    0.95950440224668454,
    0.9546774941276015,
    0.95333444514247945,
    0.95341733562653097,
    0.95091112171053427,
    0.94720344989803018,
    0.94609453502396179,
    0.99648413186838214,
    0.98936700376347253,
    0.98796849711669887,
    0.98670399251763907,
    0.97511251109787378,
    0.96209141435595047,
    0.95563437450245758,
    0.9532245701901445,
    0.94837993097111717,
    0.92691336080278564,
    0.91164413805920885,
    0.90507412088633021,
    0.89742047855477269,
    0.89574112265115413,
    0.89522110203486982
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

